$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 11 ---
# The company name for the first "cresensolutions" interview gets a trailing "1"
# to distinguish it from the new second round being added below.
$ws.Range("B11").Value = "cresensolutions - technical round1"

# --- Append new row 12 for the second cresensolutions technical round ---
$ws.Range("A12").Value = 45805
$ws.Range("B12").Value = "cresensolutions - technical round2"

$ws.Range("D12").Value = "failed"

$questions12 = "angular -  guards, subject and behaviour subject`njava - stateless and stateful operator, reflections, volatile, automatic, executer service, callable vs runnable, thread local, how to solve diamond problem of default method of interface`nspring boot - scope , transaction management, component vs bean, dependency injection, how to fix circular dependency, default scope of bean`ndouble check lock singletone pattern, circuit breaker pattern, SSE - server sent event"
$ws.Range("C12").Value = $questions12
$ws.Range("C12").WrapText = $true

# --- Keep view state consistent with the new last row ---
$ws.Range("C13").Select()
